# Update PC1 (column B) and PC2 (column C) values for rows 2-17
# on the active worksheet with newly computed results from the
# Observing FE Results on the Experimental SVM Hyperplanes script.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @{ Row = 2;  B = 0.04346061451737711;  C = -0.1949354437491366 }
    @{ Row = 3;  B = -0.02566691276721682; C = -0.345389485432891 }
    @{ Row = 4;  B = 0.09129984849215368;  C = 0.5275278924569102 }
    @{ Row = 5;  B = 0.1407055506499899;   C = -0.3602985194608619 }
    @{ Row = 6;  B = 0.1366755716739975;   C = 0.03676266634559655 }
    @{ Row = 7;  B = 0.4256831295038078;   C = 0.1351527988296268 }
    @{ Row = 8;  B = 0.3874550877375901;   C = -0.3807904616935504 }
    @{ Row = 9;  B = 0.07968010149356451;  C = 0.2667472435853089 }
    @{ Row = 10; B = 0.552950445539106;    C = -0.01248300409136403 }
    @{ Row = 11; B = 0.1573863220214999;   C = 0.06723339341275245 }
    @{ Row = 12; B = -0.3638663655288448;  C = 0.01102512915124795 }
    @{ Row = 13; B = -0.3532388052048041;  C = -0.252678093913896 }
    @{ Row = 14; B = -0.1407345591953378;  C = 0.2714236112426232 }
    @{ Row = 15; B = -0.07369506778062344; C = -0.2115150976021346 }
    @{ Row = 16; B = 0.002915719148135049; C = -0.1017092837601836 }
    @{ Row = 17; B = -0.003905432081165423;C = 0.03787105064831471 }
)

foreach ($item in $values) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
}
